# bug fix in 191
# Adds newly-run registration-history rows to the AMSIN and AMS sheets,
# and fixes the formatting/time-stamp on AMS row 33 (2024-03-08 / 189amscp)
# that had been missing its style and had a slightly wrong run-time value.

$wb = $excel.ActiveWorkbook

# Leading apostrophe forces Excel to store the value as literal text
# instead of auto-parsing date-shaped strings ("2024-03-28") into a date
# serial. Re-applying the "Normal" style afterwards clears the
# quote-prefix flag again so the cell's effective format matches the
# plain General style used by the rest of the sheet.
function SetTextCell($wsArg, $addr, $text) {
    $wsArg.Range($addr).Value = "'" + $text
    $wsArg.Range($addr).Style = "Normal"
}

function SetDateTimeCell($wsArg, $addr, $serial) {
    $wsArg.Range($addr).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsArg.Range($addr).Value = $serial
}

# ---- AMSIN sheet: append rows 43-46 ----
$wsAmsin = $wb.Worksheets.Item("AMSIN")

SetTextCell     $wsAmsin "A43" "2024-03-28"
SetDateTimeCell $wsAmsin "B43" 45379.52555458334
SetTextCell     $wsAmsin "C43" "190fstcp"
$wsAmsin.Range("D43").Value = 62
$wsAmsin.Range("E43").Value = 62
$wsAmsin.Range("F43").Value = 0
$wsAmsin.Range("G43").Value = 1.27

SetTextCell     $wsAmsin "A44" "2024-04-01"
SetDateTimeCell $wsAmsin "B44" 45383.37190483796
SetTextCell     $wsAmsin "C44" "190fnlcp"
$wsAmsin.Range("D44").Value = 62
$wsAmsin.Range("E44").Value = 62
$wsAmsin.Range("F44").Value = 0
$wsAmsin.Range("G44").Value = 1

SetTextCell     $wsAmsin "A45" "2024-05-02"
SetDateTimeCell $wsAmsin "B45" 45414.46824693287
SetTextCell     $wsAmsin "C45" "191fstcp"
$wsAmsin.Range("D45").Value = 62
$wsAmsin.Range("E45").Value = 62
$wsAmsin.Range("F45").Value = 0
$wsAmsin.Range("G45").Value = 1.13

SetTextCell     $wsAmsin "A46" "2024-05-03"
SetDateTimeCell $wsAmsin "B46" 45415.34168078704
SetTextCell     $wsAmsin "C46" "191lstcpr"
$wsAmsin.Range("D46").Value = 62
$wsAmsin.Range("E46").Value = 62
$wsAmsin.Range("F46").Value = 0
$wsAmsin.Range("G46").Value = 1.07

# ---- AMS sheet: fix row 33, append rows 34-36 ----
$wsAms = $wb.Worksheets.Item("AMS")

# Row 33 already holds 2024-03-08 / 189amscp; re-write it so it picks up
# the same style as its neighbours and carries the corrected run timestamp.
SetTextCell     $wsAms "A33" "2024-03-08"
SetDateTimeCell $wsAms "B33" 45359.72567608796
SetTextCell     $wsAms "C33" "189amscp"
$wsAms.Range("D33").Value = 62
$wsAms.Range("E33").Value = 62
$wsAms.Range("F33").Value = 0
$wsAms.Range("G33").Value = 1.1

SetTextCell     $wsAms "A34" "2024-04-01"
SetDateTimeCell $wsAms "B34" 45383.52486285879
SetTextCell     $wsAms "C34" "190betacp"
$wsAms.Range("D34").Value = 62
$wsAms.Range("E34").Value = 59
$wsAms.Range("F34").Value = 3
$wsAms.Range("G34").Value = 2.34

SetTextCell     $wsAms "A35" "2024-04-01"
SetDateTimeCell $wsAms "B35" 45383.8474333449
SetTextCell     $wsAms "C35" "190livecpp"
$wsAms.Range("D35").Value = 62
$wsAms.Range("E35").Value = 60
$wsAms.Range("F35").Value = 2
$wsAms.Range("G35").Value = 1.33

# Row 36 is appended the same way rows used to be added before this bug
# fix: plain values with no explicit style except on the run-time column.
SetTextCell     $wsAms "A36" "2024-05-03"
SetDateTimeCell $wsAms "B36" 45415.59175239856
SetTextCell     $wsAms "C36" "191betacpp"
$wsAms.Range("D36").Value = 62
$wsAms.Range("E36").Value = 62
$wsAms.Range("F36").Value = 0
$wsAms.Range("G36").Value = 1.03
